$d = $word.ActiveDocument

# This document tracks revisions by default; make sure deletions are applied
# directly to the content instead of being recorded as tracked changes.
$d.TrackRevisions = $false

# --- Change 1 -----------------------------------------------------------
# Remove the "_GoBack" bookmark from the heading paragraph at the top of
# the document (it will be re-created, empty, at the very end of the
# document further below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Change 2 -----------------------------------------------------------
# Drop the trailing "reviewer note" paragraph (the yellow-highlighted
# "[Здесь хорошо бы дописать ...]" remark) together with the blank
# paragraph that used to sit right after it, collapsing them into the
# single empty paragraph that now closes the document.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Здесь хорошо бы дописать*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Deleting the full paragraph Range (including its end-of-paragraph
    # mark) removes the paragraph outright and merges what followed it
    # into the paragraph before, leaving a single blank paragraph where
    # the two used to be.
    $target.Range.Delete()
}

# --- Re-create the "_GoBack" bookmark on the final (now empty) paragraph
# -------------------------------------------------------------------------
# The engine's Bookmarks.Add can't anchor a zero-length range inside a
# paragraph that has no runs at all, so temporarily insert a placeholder
# character to anchor the bookmark around, then remove the placeholder
# again; the (now empty) bookmark stays correctly positioned.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $last.Range
$lastRange.InsertBefore("X")

$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $last.Range
$bmRange = $d.Range($lastRange.Start, $lastRange.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $last.Range
$placeholder = $d.Range($lastRange.Start, $lastRange.Start + 1)
$placeholder.Text = ""
